$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for column C, matching style (bold/centered/bordered) of existing header cells
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Coord: normal vector scan"

# Update column B values (recomputed angles) and add column C (normal vector coords)
$ws.Range("B2").Value = 0.2198433730912588
$ws.Range("C2").Value = '[0.         0.3213351  0.94696555]'
$ws.Range("B3").Value = 1.427634881625641
$ws.Range("C3").Value = '[-0.43841952  0.52318006  0.73080158]'
$ws.Range("B4").Value = 0.08238799197386509
$ws.Range("C4").Value = '[-0.00130162  0.01543052  0.9998801 ]'
$ws.Range("B5").Value = 0.6051075703701794
$ws.Range("C5").Value = '[-0.00111763  0.2844085  -0.95870254]'
$ws.Range("B6").Value = 1.656982050048569
$ws.Range("C6").Value = '[0.72443713 0.28877145 0.62594081]'
$ws.Range("B7").Value = 0.6729724159768979
$ws.Range("C7").Value = '[-0.73447669 -0.27249856  0.62152114]'
$ws.Range("B8").Value = 0.5399908464861775
$ws.Range("C8").Value = '[0.         0.31603882 0.94874626]'
$ws.Range("B9").Value = 0.9327842266394749
$ws.Range("C9").Value = '[ 0.         -0.30952729  0.95089056]'
$ws.Range("B10").Value = 2.029716798134745
$ws.Range("C10").Value = '[-0.72718942  0.25449306  0.6375177 ]'
$ws.Range("B11").Value = 1.630453625149624
$ws.Range("C11").Value = '[ 0.73093484 -0.27747786  0.62349042]'
$ws.Range("B12").Value = 0.6597456662582706
$ws.Range("C12").Value = '[ 0.00115748 -0.28349319 -0.95897355]'
$ws.Range("B13").Value = 1.837524059595558
$ws.Range("C13").Value = '[-0.72219352 -0.28915646  0.62835106]'
$ws.Range("B14").Value = 3.635621934349079
$ws.Range("C14").Value = '[0.70292503 0.26220741 0.66116842]'
$ws.Range("B15").Value = 1.566534017546801
$ws.Range("C15").Value = '[ 0.         -0.29899075  0.95425601]'
